$d = $word.ActiveDocument

# The "DONE BY ------ ..." credits line was updated:
#   - the author's name "K.MADHUMITA" was corrected to "K.MADHUMETHA"
#   - the roll/ID number was changed from (211419106154) to (211419106153)
#   - the role tag changed from "[TEAM MEMBER 3]" to "[TEAMLEADER]"
# Apply each change as a precise, unambiguous Find & Replace so the rest
# of the paragraph / document is left untouched.

$d.Content.Find.Execute("K.MADHUMITA", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "K.MADHUMETHA", 2) | Out-Null

$d.Content.Find.Execute("(211419106154)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(211419106153)", 2) | Out-Null

$d.Content.Find.Execute("[TEAM MEMBER 3]", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[TEAMLEADER]", 2) | Out-Null
